$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for the 5 new "header" cells (AW1:BA1) by copying the
#    formatting (yellow header fill, style 1) of an existing header cell,
#    so the new cells visually match the rest of the header row.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AW1:BA1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Rename the TON / TOF header block (AC1:AL1) - wording only, same columns
# ---------------------------------------------------------------------------
$ws.Range("AC1").Value = "Turnover number(CO) <TON CO>"
$ws.Range("AD1").Value = "Turnover number(CH4) <TON CH4>"
$ws.Range("AE1").Value = "Turnover number(H2) <TON H2>"
$ws.Range("AF1").Value = "Turnover number(HCOOH) <TON HCOOH>"
$ws.Range("AG1").Value = "Turnover number(MeOH) <TON MeOH>"
$ws.Range("AH1").Value = "Turnover frequency(CO) <TOF CO>"
$ws.Range("AI1").Value = "Turnover frequency(CH4) <TOF CH4>"
$ws.Range("AJ1").Value = "Turnover frequency(H2) <TOF H2>"
$ws.Range("AK1").Value = "Turnover frequency(HCOOH) <TOF HCOOH>"
$ws.Range("AL1").Value = "Turnover frequency(MeOH) <TOF MeOH>"

# ---------------------------------------------------------------------------
# 3. The Quantum-yield block grows from 5 to 10 columns (an "_inchikey"
#    column followed by a new "_molfile" column for each gas), and the
#    trailing columns (Quantum yield total / Selectivity / Details /
#    Included / BasePageName) move 5 columns to the right, into the cells
#    whose formatting was prepared in step 1. Simply overwrite every header
#    cell from AM1 to BA1 with its final text.
# ---------------------------------------------------------------------------
$ws.Range("AM1").Value = "Quantum yield(CO)_inchikey <Φ CO>"
$ws.Range("AN1").Value = "Quantum yield at (CO)_molfile"
$ws.Range("AO1").Value = "Quantum yield(CH4)_inchikey <Φ CH4>"
$ws.Range("AP1").Value = "Quantum yield at (CH4)_molfile"
$ws.Range("AQ1").Value = "Quantum yield(H2)_inchikey <Φ H2>"
$ws.Range("AR1").Value = "Quantum yield at (H2)_molfile"
$ws.Range("AS1").Value = "Quantum yield(HCOOH)_inchikey <Φ HCOOH>"
$ws.Range("AT1").Value = "Quantum yield at (HCOOH)_molfile"
$ws.Range("AU1").Value = "Quantum yield(MeOH)_inchikey <Φ MeOH>"
$ws.Range("AV1").Value = "Quantum yield at (MeOH)_molfile"
$ws.Range("AW1").Value = "Quantum yield total <Φ all>"
$ws.Range("AX1").Value = "Selectivity <auto-generated-Selectivity>"
$ws.Range("AY1").Value = "Details <details>"
$ws.Range("AZ1").Value = "Included <include>"
$ws.Range("BA1").Value = "BasePageName <BasePageName>"

# ---------------------------------------------------------------------------
# 4. Column width tweaks (A, B, and the new AK "Turnover frequency(HCOOH)")
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 33.166666666666664
$ws.Columns("B").ColumnWidth = 254.8
$ws.Columns("AK").ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# 5. Selection tweak
# ---------------------------------------------------------------------------
$ws.Range("A2:XFD3").Select() | Out-Null
